# Adicionado cantidad de engranajes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: gear tooth counts (Z), Column B: cantidad (quantity)
$data = @(
    @(21, 1),
    @(24, 2),
    @(25, 1),
    @(28, 1),
    @(30, 1),
    @(32, 3),
    @(35, 1),
    @(36, 1),
    @(39, 1),
    @(40, 2),
    @(44, 1),
    @(45, 1),
    @(47, 1),
    @(48, 1),
    @(50, 1),
    @(51, 1),
    @(52, 1),
    @(53, 1),
    @(55, 1),
    @(56, 1),
    @(58, 1),
    @(60, 2),
    @(61, 1),
    @(62, 1),
    @(64, 1),
    @(66, 1),
    @(70, 1),
    @(71, 1),
    @(72, 1),
    @(73, 1),
    @(74, 1),
    @(75, 1),
    @(76, 1),
    @(79, 1),
    @(80, 1),
    @(84, 1),
    @(86, 1),
    @(88, 1),
    @(90, 1),
    @(96, 1),
    @(100, 1)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Update selection to match the target state
$ws.Range("B25").Select()
